$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Row 9 corrections ("Metodos Norma 1, 2, inf en vector"): fill in missing LoC real (C9) and Hora de Fin (F9)
$ws.Range("C9").Value = 24
$ws.Range("F9").Value = 0.73611111111111116

# Widen column G (it's no longer best-fit, just a manually set width)
$ws.Columns.Item(7).ColumnWidth = 11.42578125

# Update the active selection to reflect where the editor ended up working
$ws.Range("H9").Select()

$wb.Save()
